$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(4, 4, 2, 6, 21, 2)
$row = 5
foreach ($v in $values) {
    $ws.Cells.Item($row, 1).Value = $v
    $row++
}

$ws.Range("A10").Select()
